$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 93 (shifts existing rows 93..153 down to 94..154)
$ws.Rows("93:93").Insert()

# Populate the new row 93 with the new weekly price-report entry
$ws.Cells.Item(93, 1).Value = 4
$ws.Cells.Item(93, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value = "Los Lagos"
$ws.Cells.Item(93, 4).Value = "05/22/2023"
$ws.Cells.Item(93, 5).Value = 10
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100104
$ws.Cells.Item(93, 8).Value = "Frutos de pepita"
$ws.Cells.Item(93, 9).Value = 100104003
$ws.Cells.Item(93, 10).Value = "Membrillo"
$ws.Cells.Item(93, 11).Value = "Champion"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 200
$ws.Cells.Item(93, 14).Value = 13000
$ws.Cells.Item(93, 15).Value = 14000
$ws.Cells.Item(93, 16).Value = 13500
$ws.Cells.Item(93, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(93, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(93, 19).Value = 750
$ws.Cells.Item(93, 20).Value = 18
